$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows 8 and 9 (Resolving-Mac target-cluster rows no longer present)
$ws.Rows("8:9").Delete()

# Update remaining data rows (2-7) with new TPM-derived values
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Bdnf"
$ws.Range("C2").Value = "Ntrk2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.020961333333334
$ws.Range("H2").Value = 6.062884
$ws.Range("I2").Value = 0.2503572190582515
$ws.Range("J2").Value = 0.2503572190582515
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.8220243333333334
$ws.Range("N2").Value = 2.466073
$ws.Range("O2").Value = 0.03815249372618141
$ws.Range("P2").Value = 0.03815249372618141
$ws.Range("Q2").Value = 1.661279392725778
$ws.Range("R2").Value = 14.951514534532
$ws.Range("S2").Value = 0.009551752229424166
$ws.Range("T2").Value = 0.009551752229424166
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Bdnf"
$ws.Range("C3").Value = "Ntrk2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.020961333333334
$ws.Range("H3").Value = 6.062884
$ws.Range("I3").Value = 0.2503572190582515
$ws.Range("J3").Value = 0.2503572190582515
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 14.52590566666666
$ws.Range("N3").Value = 43.57771699999999
$ws.Range("O3").Value = 0.6741887099221348
$ws.Range("P3").Value = 0.6741887099221348
$ws.Range("Q3").Value = 29.35629368398089
$ws.Range("R3").Value = 264.206643155828
$ws.Range("S3").Value = 0.1687880105365759
$ws.Range("T3").Value = 0.1687880105365759
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Bdnf"
$ws.Range("C4").Value = "Ntrk2"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.020961333333334
$ws.Range("H4").Value = 6.062884
$ws.Range("I4").Value = 0.2503572190582515
$ws.Range("J4").Value = 0.2503572190582515
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 6.197826333333334
$ws.Range("N4").Value = 18.593479
$ws.Range("O4").Value = 0.2876587963516838
$ws.Range("P4").Value = 0.2876587963516837
$ws.Range("Q4").Value = 12.52556737038178
$ws.Range("R4").Value = 112.730106333436
$ws.Range("S4").Value = 0.07201745629225147
$ws.Range("T4").Value = 0.07201745629225145
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Bdnf"
$ws.Range("C5").Value = "Ntrk2"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 6.051349666666667
$ws.Range("H5").Value = 18.154049
$ws.Range("I5").Value = 0.7496427809417484
$ws.Range("J5").Value = 0.7496427809417485
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.8220243333333334
$ws.Range("N5").Value = 2.466073
$ws.Range("O5").Value = 0.03815249372618141
$ws.Range("P5").Value = 0.03815249372618141
$ws.Range("Q5").Value = 4.974356675508556
$ws.Range("R5").Value = 44.76921007957701
$ws.Range("S5").Value = 0.02860074149675724
$ws.Range("T5").Value = 0.02860074149675725
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Bdnf"
$ws.Range("C6").Value = "Ntrk2"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 6.051349666666667
$ws.Range("H6").Value = 18.154049
$ws.Range("I6").Value = 0.7496427809417484
$ws.Range("J6").Value = 0.7496427809417485
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 14.52590566666666
$ws.Range("N6").Value = 43.57771699999999
$ws.Range("O6").Value = 0.6741887099221348
$ws.Range("P6").Value = 0.6741887099221348
$ws.Range("Q6").Value = 87.90133441401477
$ws.Range("R6").Value = 791.1120097261329
$ws.Range("S6").Value = 0.5054006993855589
$ws.Range("T6").Value = 0.5054006993855589
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Bdnf"
$ws.Range("C7").Value = "Ntrk2"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 6.051349666666667
$ws.Range("H7").Value = 18.154049
$ws.Range("I7").Value = 0.7496427809417484
$ws.Range("J7").Value = 0.7496427809417485
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 6.197826333333334
$ws.Range("N7").Value = 18.593479
$ws.Range("O7").Value = 0.2876587963516838
$ws.Range("P7").Value = 0.2876587963516837
$ws.Range("Q7").Value = 37.50521431627456
$ws.Range("R7").Value = 337.5469288464711
$ws.Range("S7").Value = 0.2156413400594323
$ws.Range("T7").Value = 0.2156413400594323
